$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-10 Sunday" "2023-12-11 Monday"

Replace-Text "28×96=" "97×49="
Replace-Text "26×34=" "79×60="
Replace-Text "32×69=" "69×70="
Replace-Text "21×16=" "28×29="
Replace-Text "36×22=" "31×49="

Replace-Text "13×40=" "21×93="
Replace-Text "46×67=" "58×62="
Replace-Text "51×57=" "51×84="
Replace-Text "98×18=" "23×23="
Replace-Text "73×43=" "82×55="

Replace-Text "62×27=" "43×66="
Replace-Text "26×63=" "17×47="
Replace-Text "64×66=" "85×16="
Replace-Text "66×52=" "25×55="
Replace-Text "99×61=" "11×51="

Replace-Text "37×90=" "38×20="
Replace-Text "57×53=" "70×31="
Replace-Text "62×80=" "27×90="
Replace-Text "78×67=" "84×19="
Replace-Text "38×37=" "70×52="

Replace-Text "12×94=" "14×25="
Replace-Text "93×44=" "73×30="
Replace-Text "43×53=" "33×91="
Replace-Text "69×93=" "19×83="
Replace-Text "93×51=" "75×26="
